$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.298.17"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.601.90"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "1.828.62"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "1.598.44"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "26.317.34"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.61%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "1.449.31"
$ws.Range("E33").Value = "  +7.88%  "
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.564"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.926"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").Value = "1.740.74"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.757"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -3.64%  "
